$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Before:
#   1  Week 6: DOM Manipulation in JavaScript
#   2  (empty paragraph containing the _GoBack bookmark)
#   3  Lab 1: DOM Selection and Manipulation
#   4  Objective: Practice selecting and manipulating DOM elements.
#   5  Steps:
#   6    1. Select an element by `id` and change its content.
#   7    2. Change the style of an element.
#   8  (empty paragraph)
#   9  (paragraph containing the lab-1 screenshot drawing)
#  10  Lab 2: Event Handling
#       ...
#
# After:
#   1  Week 6: DOM Manipulation in JavaScript
#   2  (new empty paragraph)
#   3  Lab 1: DOM Selection and Manipulation
#   4  Objective: Practice selecting and manipulating DOM elements.
#   5  Steps:
#   6    1. Select an element by `id` and change its content.
#   7    2. Change the style of an element.
#   8  (new empty paragraph)
#   9  (empty paragraph containing the _GoBack bookmark)
#  10  (paragraph containing the lab-1 screenshot drawing)
#  11  (new empty paragraph)
#  12  (new empty paragraph)
#  13  Lab 2: Event Handling
#       ...
#
# i.e. the "Lab 1" block (5 paragraphs) moves from after the bookmark
# paragraph to before it, a new blank paragraph is added right after the
# "Week 6" title, another blank paragraph separates the moved block from the
# bookmark paragraph (replacing the blank paragraph that used to sit there),
# and two new blank paragraphs are inserted between the lab-1 drawing and
# "Lab 2: Event Handling".
# ---------------------------------------------------------------------------

# Step 1: cut the "Lab 1" block (Lab1 title, Objective, Steps, step 1, step 2)
# and move it to sit right after the "Week 6" heading, i.e. before the
# bookmark paragraph.
$blockStart = $d.Paragraphs.Item(3)
$blockEnd   = $d.Paragraphs.Item(7)
$blockRange = $d.Range($blockStart.Range.Start, $blockEnd.Range.End)
$blockRange.Cut() | Out-Null

$titlePara = $d.Paragraphs.Item(1)
$pasteTarget = $d.Range($titlePara.Range.End, $titlePara.Range.End)
$pasteTarget.Paste() | Out-Null

# Now: 1 Week6, 2 Lab1, 3 Objective, 4 Steps, 5 step1, 6 step2,
#      7 bookmark, 8 old-empty, 9 drawing, 10 Lab2, ...

# Step 2: add a new blank paragraph right after "Week 6" (before "Lab 1").
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

# Now: 1 Week6, 2 new-empty, 3 Lab1, 4 Objective, 5 Steps, 6 step1, 7 step2,
#      8 bookmark, 9 old-empty, 10 drawing, 11 Lab2, ...

# Step 3: add a new blank paragraph right after the moved block (before the
# bookmark paragraph).
$step2Para = $d.Paragraphs.Item(7)
$step2Para.Range.InsertParagraphAfter() | Out-Null

# Now: 1 Week6, 2 new-empty, 3 Lab1, 4 Objective, 5 Steps, 6 step1, 7 step2,
#      8 new-empty, 9 bookmark, 10 old-empty, 11 drawing, 12 Lab2, ...

# Step 4: remove the old blank paragraph that used to separate the Lab 1
# block from its screenshot (it sat right before the drawing, and now sits
# between the bookmark paragraph and the drawing).
$oldEmpty = $d.Paragraphs.Item(10)
$oldEmpty.Range.Delete() | Out-Null

# Now: 1 Week6, 2 new-empty, 3 Lab1, 4 Objective, 5 Steps, 6 step1, 7 step2,
#      8 new-empty, 9 bookmark, 10 drawing, 11 Lab2, ...

# Step 5: insert two new blank paragraphs right after the drawing paragraph
# (before "Lab 2: Event Handling").
$drawingPara = $d.Paragraphs.Item(10)
$drawingPara.Range.InsertParagraphAfter() | Out-Null
$drawingPara = $d.Paragraphs.Item(10)
$drawingPara.Range.InsertParagraphAfter() | Out-Null

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
